$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of labels in column A (row 1 is header "Counts" in B1, unchanged).
# Column B (counts) remain unchanged; only the text labels are being re-ordered
# to match the re-generated shared string table order.
$values = @(
    'мелочный товар',
    'шелковый товар',
    'съестной припасы',
    'бумажный товар',
    'деревянный товар',
    'крестьянский товар',
    'москательный товар',
    'мелкий товар',
    'лавочный товар',
    'рукоделие',
    'шерстяной товар',
    'красный товар',
    'гарусный товар',
    'мелочь',
    'особливый товар',
    'деревенский товар',
    'серебреный товар',
    'крамными товар',
    'небогатый товар',
    'железный товар',
    'мясо',
    'приуготовлять',
    'пушной товар',
    'щепетильный товар',
    'нужный товар',
    'набойчатый товар',
    'суровский товар',
    'недорогой товар',
    'внутренний товар',
    'питейный припасы',
    'медный товар',
    'привозный товар',
    'оловянный товар',
    'произрастание',
    'заморский товар',
    'купецкий товар',
    'галантерейный товар',
    'надлежащий товар',
    'домовый товар',
    'харчевой припасы',
    'рукодельный товар',
    'меховой товар',
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
